# Auto-generated Excel COM-interop script
# Applies leve-profit market-price updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# ALC: set 48 cell(s)
$ALC_values = @{
    "H11" = 7588.5713
    "I11" = 7588.5713
    "K11" = 7588.5713
    "M11" = -7448.5713
    "H18" = 1220.2
    "I18" = 1333.6666
    "J18" = 1050
    "K18" = 1333.6666
    "L18" = 1050
    "M18" = -1049.6666
    "N18" = -1618
    "H64" = 66750
    "I64" = 5000
    "J64" = 169666.67
    "K64" = 5000
    "L64" = 169666.67
    "M64" = -4752
    "N64" = -170162.67
    "H67" = 66750
    "I67" = 5000
    "J67" = 169666.67
    "K67" = 5000
    "L67" = 169666.67
    "M67" = -4142
    "N67" = -171382.67
    "H113" = 1571.4286
    "I113" = 0
    "J113" = 1571.4286
    "K113" = 0
    "N113" = -8079.4286
    "H116" = 3476.7778
    "I116" = 2541.5715
    "J116" = 6750
    "K116" = 2541.5715
    "L116" = 6750
    "M116" = 900.4285
    "N116" = -13634
    "H132" = 246274.81
    "I132" = 288322.06
    "K132" = 864966.1799999999
    "M132" = -862436.1799999999
    "H138" = 2086169.5
    "I138" = 3126999
    "J138" = 4510.125
    "K138" = 9380997
    "L138" = 13530.375
    "M138" = -9375857
    "N138" = -23810.375
}
foreach ($cell in $ALC_values.Keys) {
    $ws.Range($cell).Value = $ALC_values[$cell]
}
# ALC: clear 1 cell(s)
foreach ($cell in @("M113")) {
    $ws.Range($cell).ClearContents()
}

$ws = $wb.Worksheets.Item("ARM")

# ARM: set 32 cell(s)
$ARM_values = @{
    "H55" = 17850.334
    "J55" = 17850.334
    "L55" = 17850.334
    "N55" = -18480.334
    "H61" = 1415.1034
    "I61" = 1426.2
    "J61" = 1390.4445
    "K61" = 1426.2
    "L61" = 1390.4445
    "M61" = -1214.2
    "N61" = -1814.4445
    "H97" = 688.75
    "I97" = 542
    "J97" = 933.3333
    "K97" = 542
    "L97" = 933.3333
    "M97" = -46
    "N97" = -1925.3333
    "H132" = 2083.9268
    "I132" = 1942.8611
    "J132" = 3099.6
    "K132" = 5828.5833
    "L132" = 9298.799999999999
    "M132" = -3298.5833
    "N132" = -14358.8
    "H136" = 1415.1034
    "I136" = 1426.2
    "J136" = 1390.4445
    "K136" = 4278.6
    "L136" = 4171.333500000001
    "M136" = -1728.6
    "N136" = -9271.333500000001
}
foreach ($cell in $ARM_values.Keys) {
    $ws.Range($cell).Value = $ARM_values[$cell]
}

$ws = $wb.Worksheets.Item("BSM")

# BSM: set 29 cell(s)
$BSM_values = @{
    "H80" = 178.23529
    "I80" = 63
    "J80" = 213.6923
    "K80" = 63
    "L80" = 213.6923
    "M80" = 935
    "N80" = -2209.6923
    "H83" = 178.23529
    "I83" = 63
    "J83" = 213.6923
    "K83" = 315
    "L83" = 1068.4615
    "M83" = 4677
    "N83" = -11052.4615
    "H94" = 623.2353000000001
    "I94" = 491.85715
    "K94" = 491.85715
    "M94" = -40.85714999999999
    "H99" = 1007
    "I99" = 968.4286
    "J99" = 1074.5
    "K99" = 968.4286
    "L99" = 1074.5
    "M99" = 529.5714
    "N99" = -4070.5
    "H134" = 54180.05
    "J134" = 1756.25
    "L134" = 5268.75
    "N134" = -10338.75
}
foreach ($cell in $BSM_values.Keys) {
    $ws.Range($cell).Value = $BSM_values[$cell]
}

$ws = $wb.Worksheets.Item("CRP")

# CRP: set 81 cell(s)
$CRP_values = @{
    "H22" = 275.29166
    "I22" = 184.36842
    "J22" = 620.8
    "K22" = 184.36842
    "L22" = 620.8
    "M22" = 165.63158
    "N22" = -1320.8
    "H31" = 2235
    "I31" = 2064
    "J31" = 2357.1428
    "K31" = 2064
    "L31" = 2357.1428
    "M31" = -1769
    "N31" = -2947.1428
    "H34" = 2235
    "I34" = 2064
    "J34" = 2357.1428
    "K34" = 2064
    "L34" = 2357.1428
    "M34" = -1862
    "N34" = -2761.1428
    "H41" = 19600
    "J41" = 19600
    "L41" = 19600
    "N41" = -20456
    "H50" = 10175.375
    "J50" = 10343.286
    "L50" = 10343.286
    "N50" = -11593.286
    "H51" = 10400.143
    "I51" = 9100
    "K51" = 9100
    "M51" = -8364
    "H59" = 16562.125
    "J59" = 16562.125
    "L59" = 16562.125
    "N59" = -18852.125
    "H60" = 9345.111000000001
    "I60" = 6000
    "J60" = 10300.857
    "K60" = 6000
    "L60" = 10300.857
    "M60" = -5489
    "N60" = -11322.857
    "H61" = 10400.143
    "I61" = 9100
    "K61" = 9100
    "M61" = -8752
    "H74" = 0
    "J74" = 0
    "H77" = 0
    "J77" = 0
    "H99" = 2275
    "I99" = 2066.6667
    "J99" = 2900
    "K99" = 2066.6667
    "L99" = 2900
    "M99" = -568.6667000000002
    "N99" = -5896
    "H109" = 18285
    "J109" = 18285
    "L109" = 18285
    "N109" = -20365
    "H126" = 2275
    "I126" = 2066.6667
    "J126" = 2900
    "K126" = 6200.000100000001
    "L126" = 8700
    "M126" = -3730.000100000001
    "N126" = -13640
    "H132" = 2496.8262
    "I132" = 2337.6365
    "J132" = 5999
    "K132" = 7012.9095
    "L132" = 17997
    "M132" = -4482.9095
    "N132" = -23057
    "H134" = 5057.2
    "I134" = 5057.2
    "K134" = 15171.6
    "M134" = -12636.6
}
foreach ($cell in $CRP_values.Keys) {
    $ws.Range($cell).Value = $CRP_values[$cell]
}
# CRP: clear 2 cell(s)
foreach ($cell in @("N74", "N77")) {
    $ws.Range($cell).ClearContents()
}

$ws = $wb.Worksheets.Item("CUL")

# CUL: set 15 cell(s)
$CUL_values = @{
    "H5" = 1004.34784
    "I5" = 1033.3334
    "K5" = 3100.0002
    "M5" = -2988.0002
    "H135" = 1004.34784
    "I135" = 1033.3334
    "K135" = 9300.000599999999
    "M135" = -6765.000599999999
    "H140" = 1506.25
    "I140" = 978.125
    "J140" = 2562.5
    "K140" = 2934.375
    "L140" = 7687.5
    "M140" = 2245.625
    "N140" = -18047.5
}
foreach ($cell in $CUL_values.Keys) {
    $ws.Range($cell).Value = $CUL_values[$cell]
}

$ws = $wb.Worksheets.Item("GSM")

# GSM: set 7 cell(s)
$GSM_values = @{
    "H132" = 2591.9
    "I132" = 2233.75
    "J132" = 4024.5
    "K132" = 6701.25
    "L132" = 12073.5
    "M132" = -4171.25
    "N132" = -17133.5
}
foreach ($cell in $GSM_values.Keys) {
    $ws.Range($cell).Value = $GSM_values[$cell]
}

$ws = $wb.Worksheets.Item("LTW")

# LTW: set 22 cell(s)
$LTW_values = @{
    "H22" = 517.2857
    "J22" = 516.3333
    "L22" = 516.3333
    "N22" = -1106.3333
    "H27" = 517.2857
    "J27" = 516.3333
    "L27" = 516.3333
    "N27" = -730.3333
    "H46" = 1747.5264
    "I46" = 1183.4166
    "J46" = 2714.5715
    "K46" = 1183.4166
    "L46" = 2714.5715
    "M46" = -995.4166
    "N46" = -3090.5715
    "H136" = 3958.125
    "I136" = 1790
    "J136" = 5644.4443
    "K136" = 5370
    "L136" = 16933.3329
    "M136" = -2820
    "N136" = -22033.3329
}
foreach ($cell in $LTW_values.Keys) {
    $ws.Range($cell).Value = $LTW_values[$cell]
}

$ws = $wb.Worksheets.Item("WVR")

# WVR: set 6 cell(s)
$WVR_values = @{
    "H136" = 16866.385
    "I136" = 41456.8
    "J136" = 1497.375
    "K136" = 124370.4
    "L136" = 4492.125
    "N136" = -9592.125
}
foreach ($cell in $WVR_values.Keys) {
    $ws.Range($cell).Value = $WVR_values[$cell]
}
